$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 160 ("Vega Modelo de Temuco" /
# "Pepino dulce" weekly price sheet). This shifts the existing rows 160-165
# down to 161-166, preserving all of their original data untouched.
$ws.Rows(160).Insert()

# Populate the newly-inserted row 160 with this week's price record. Most
# columns repeat the constant values used throughout this sub-sheet; only the
# date (D) and the volume/price columns (J, K, L, M, P) carry new data.
$ws.Cells.Item(160, 1).Value2  = 10
$ws.Cells.Item(160, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(160, 3).Value2  = 'La Araucanía'
$ws.Cells.Item(160, 4).Value2  = 44568
$ws.Cells.Item(160, 5).Value2  = 9
$ws.Cells.Item(160, 6).Value2  = 100112043
$ws.Cells.Item(160, 7).Value2  = 'Pepino dulce'
$ws.Cells.Item(160, 8).Value2  = 'Cultivar IV Región'
$ws.Cells.Item(160, 9).Value2  = 'Primera'
$ws.Cells.Item(160, 10).Value2 = 100
$ws.Cells.Item(160, 11).Value2 = 25000
$ws.Cells.Item(160, 12).Value2 = 25000
$ws.Cells.Item(160, 13).Value2 = 25000
$ws.Cells.Item(160, 14).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(160, 15).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(160, 16).Value2 = 1389
$ws.Cells.Item(160, 17).Value2 = 18
$ws.Cells.Item(160, 18).Value2 = 'Hortaliza'

# Match the date-formatted style already used by column D on this sheet.
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(161, 4).NumberFormat
